$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tuesday")
$ws.Activate()
$excel.ActiveWindow.Zoom = 64

$c16 = $ws.Range("C16").Value2
$e16 = $ws.Range("E16").Value2
$g16 = $ws.Range("G16").Value2

$ws.Range("C11").Value = $c16
$ws.Range("E11").Value = $e16
$ws.Range("G11").Value = $g16

$ws.Range("C16").Value = ""
$ws.Range("E16").Value = ""
$ws.Range("G16").Value = ""
